# Gate_Closure_Trigger.xlsx — update gate-closure trigger values and
# restore the saved cursor/selection position.
# (Adding timeout option of 20 seconds when attempting to download xml
#  files from rivergages.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# IHNC Surge Barrier trigger: 3 -> 5
$ws.Range("C3").Value = 5

# Bayou Dupre Sector Gate trigger: 0 -> 2
$ws.Range("C7").Value = 2

# Empire Lock trigger: 0 -> 2
$ws.Range("C10").Value = 2

# Leave the selection where the editor left it when saving (C11)
$ws.Range("C11").Select() | Out-Null
